$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Version bump
$meta.Range("B3").Value = "6.0.0"

# Publication date refresh
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was previously blank
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was a duplicate "Contact / No display for ContactDetail" row.
# Turn it into the new "Jurisdiction / United States of America" row...
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# ...and remove the other duplicate "Contact" row (row 11), which shifts
# all following rows (Description, Purpose, Copyright, ...) up by one.
$meta.Rows.Item(11).Delete()

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Root element's Short/Definition now reflect the specific extension
# instead of the generic "Extension" / "An Extension" placeholder text.
$elements.Range("K2").Value = "Element Source Classification"
$elements.Range("L2").Value = "Classification of the origin of the data value associated to a given element in a FHIR resource. The intent is that this extension value should be populated with a code from the process-meta-source-classification valueset."
